$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) cells stay text even for numeric-looking values like "1.031"
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 4).Value = '27.468.37'
$ws.Cells.Item(2, 5).Value = '  +4.24%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.836.69'
$ws.Cells.Item(3, 5).Value = '  +3.52%  '

# Row 4
$ws.Cells.Item(4, 4).Value = '1.031'
$ws.Cells.Item(4, 5).Value = '  +2.86%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '318.24'
$ws.Cells.Item(5, 5).Value = '  +3.22%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '1.027'
$ws.Cells.Item(6, 5).Value = '  +2.48%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '0.4369'
$ws.Cells.Item(7, 5).Value = '  +3.27%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '0.3724'
$ws.Cells.Item(8, 5).Value = '  +3.24%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '0.07365'
$ws.Cells.Item(9, 5).Value = '  +3.34%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '0.8736'
$ws.Cells.Item(10, 5).Value = '  +4.34%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '21.42'
$ws.Cells.Item(11, 5).Value = '  +5.17%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '1.863.55'
$ws.Cells.Item(12, 5).Value = '  +5.11%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '5.479'
$ws.Cells.Item(13, 5).Value = '  +4.42%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '6.671'
$ws.Cells.Item(14, 5).Value = '  +3.50%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '0.07133'
$ws.Cells.Item(15, 5).Value = '  +3.45%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '82.55'
$ws.Cells.Item(16, 5).Value = '  +4.52%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '1.030'
$ws.Cells.Item(17, 5).Value = '  +2.74%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '0.000008987'
$ws.Cells.Item(18, 5).Value = '  +3.89%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '1.023'
$ws.Cells.Item(19, 5).Value = '  +2.14%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '15.38'
$ws.Cells.Item(20, 5).Value = '  +3.21%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '27.526.46'
$ws.Cells.Item(21, 5).Value = '  +4.38%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '5.228'
$ws.Cells.Item(22, 5).Value = '  +2.56%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '11.18'
$ws.Cells.Item(23, 5).Value = '  +2.20%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '2.074.72'
$ws.Cells.Item(24, 5).Value = '  +4.24%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '156.89'
$ws.Cells.Item(25, 5).Value = '  +3.32%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +6.79%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '18.70'
$ws.Cells.Item(27, 5).Value = '  +3.95%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '5.243'
$ws.Cells.Item(28, 5).Value = '  +3.55%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '1.932'
$ws.Cells.Item(29, 5).Value = '  +4.93%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '116.23'
$ws.Cells.Item(30, 5).Value = '  +1.73%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '0.09055'
$ws.Cells.Item(31, 5).Value = '  +2.53%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '1.207'
$ws.Cells.Item(32, 5).Value = '  +7.73%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '0.7633'
$ws.Cells.Item(33, 5).Value = '  +5.14%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '4.483'
$ws.Cells.Item(34, 5).Value = '  +3.97%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '2.875'
$ws.Cells.Item(35, 5).Value = '  +5.03%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  +2.60%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '1.146'
$ws.Cells.Item(37, 5).Value = '  +5.90%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '0.01966'
$ws.Cells.Item(38, 5).Value = '  +4.31%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '0.05245'
$ws.Cells.Item(39, 5).Value = '  +2.74%  '

# Row 40
$ws.Cells.Item(40, 2).Value = 'MXToken'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(40, 4).Value = '2.812'
$ws.Cells.Item(40, 5).Value = '  +8.18%  '

# Row 41
$ws.Cells.Item(41, 2).Value = 'TheSandbox'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(41, 4).Value = '0.5166'
$ws.Cells.Item(41, 5).Value = '  +5.21%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '0.1665'
$ws.Cells.Item(42, 5).Value = '  +3.53%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '6.607'
$ws.Cells.Item(43, 5).Value = '  +4.37%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '8.504'
$ws.Cells.Item(44, 5).Value = '  +5.81%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '108.92'
$ws.Cells.Item(45, 5).Value = '  +4.07%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '10.55'
$ws.Cells.Item(46, 5).Value = '  +3.83%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '1.029'
$ws.Cells.Item(47, 5).Value = '  +2.73%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '1.699'
$ws.Cells.Item(48, 5).Value = '  +4.70%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '0.4635'
$ws.Cells.Item(49, 5).Value = '  +4.41%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '1.902'
$ws.Cells.Item(50, 5).Value = '  +12.17%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '0.06347'
$ws.Cells.Item(51, 5).Value = '  +2.84%  '

# Reset Price column style so no stray text-format style lingers on the cells
$ws.Range("D2:D51").Style = "Normal"
